# Update "想去人数" (interest count) values in the F column across the
# 展览 (Exhibitions), 演出 (Shows) and 全部类型 (All Types) sheets.

$wb = $excel.ActiveWorkbook

# 展览 sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1128
$ws1.Range("F5").Value = 183
$ws1.Range("F7").Value = 234
$ws1.Range("F11").Value = 518
$ws1.Range("F12").Value = 542
$ws1.Range("F14").Value = 12822
$ws1.Range("F15").Value = 8
$ws1.Range("F16").Value = 5265

# 演出 sheet
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 4
$ws2.Range("F3").Value = 108

# 全部类型 sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 4
$ws4.Range("F5").Value = 1128
$ws4.Range("F6").Value = 183
$ws4.Range("F8").Value = 234
$ws4.Range("F12").Value = 518
$ws4.Range("F13").Value = 542
$ws4.Range("F15").Value = 12822
$ws4.Range("F16").Value = 108
$ws4.Range("F18").Value = 8
$ws4.Range("F19").Value = 5265
